# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Sat Oct  5 16:51:59 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.109.84'
$ws.Range('E2').Value = '  -0.18%  '

# Row 3
$ws.Range('E3').Value = '  -0.51%  '

# Row 4
$ws.Range('E4').Value = '  -0.13%  '

# Row 5
$ws.Range('D5').Value = '''562.20'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.25%  '

# Row 6
$ws.Range('D6').Value = '''142.40'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.83%  '

# Row 7
$ws.Range('E7').Value = '  +0.12%  '

# Row 8
$ws.Range('E8').Value = '  -0.53%  '

# Row 9
$ws.Range('E9').Value = '  -0.07%  '

# Row 10
$ws.Range('E10').Value = '  -2.09%  '

# Row 11
$ws.Range('D11').Value = '''5.30'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.93%  '

# Row 12
$ws.Range('E12').Value = '  -1.55%  '

# Row 13
$ws.Range('D13').Value = '''25.54'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.21%  '

# Row 14
$ws.Range('D14').Value = '''0.0000174'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.02%  '

# Row 15
$ws.Range('D15').Value = '2.846.50'
$ws.Range('E15').Value = '  -0.58%  '

# Row 16
$ws.Range('D16').Value = '62.169.41'
$ws.Range('E16').Value = '  +0.33%  '

# Row 17
$ws.Range('D17').Value = '2.406.89'
$ws.Range('E17').Value = '  -0.91%  '

# Row 18
$ws.Range('E18').Value = '  +0.53%  '

# Row 19
$ws.Range('E19').Value = '  -1.37%  '

# Row 20
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '''6.83'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.50%  '

# Row 21
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '''320.54'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.44%  '

# Row 22
$ws.Range('E22').Value = '  -0.10%  '

# Row 23
$ws.Range('D23').Value = '''66.02'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.65%  '

# Row 24
$ws.Range('E24').Value = '  -1.98%  '

# Row 25
$ws.Range('D25').Value = '''8.84'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.17%  '

# Row 26
$ws.Range('D26').Value = '''569.43'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.04%  '

# Row 27
$ws.Range('E27').Value = '  +0.47%  '

# Row 28
$ws.Range('D28').Value = '2.528.05'

# Row 29
$ws.Range('D29').Value = '0.0₃0939'
$ws.Range('E29').Value = '  +0.37%  '

# Row 30
$ws.Range('D30').Value = '''8.16'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.02%  '

# Row 31
$ws.Range('E31').Value = '  -3.29%  '

# Row 32
$ws.Range('E32').Value = '  -0.57%  '

# Row 33
$ws.Range('D33').Value = '''1.86'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.28%  '

# Row 34
$ws.Range('E34').Value = '  -2.62%  '

# Row 35
$ws.Range('E35').Value = '  +0.24%  '

# Row 36
$ws.Range('E36').Value = '  -2.87%  '

# Row 37
$ws.Range('E37').Value = '  -5.97%  '

# Row 38
$ws.Range('E38').Value = '  -1.24%  '

# Row 39
$ws.Range('D39').Value = '''151.70'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.85%  '

# Row 40
$ws.Range('D40').Value = '''18.60'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.24%  '

# Row 41
$ws.Range('E41').Value = '  -9.56%  '

# Row 42
$ws.Range('D42').Value = '''0.993'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.62%  '

# Row 43
$ws.Range('D43').Value = '''2.27'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.02%  '

# Row 44
$ws.Range('D44').Value = '''147.73'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.24%  '

# Row 46
$ws.Range('D46').Value = '''0.0532'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.78%  '

# Row 47
$ws.Range('E47').Value = '  -2.86%  '

# Row 48
$ws.Range('E48').Value = '  -0.27%  '

# Row 49
$ws.Range('D49').Value = '''0.0915'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.21%  '

# Row 50
$ws.Range('E50').Value = '  -0.85%  '

# Row 51
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '''11.53'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.13%  '

